$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the statistic headers to their new (de-duplicated) labels
$ws.Range("B1").Value = "New Zealand_priceprice"
$ws.Range("D1").Value = "New Zealand_pointspoints"

# Remove the spacer column C (was the stray "_1" label with no data),
# which shifts the "points" column (D) left into C
$ws.Range("C1").EntireColumn.Delete()

# Match the new column widths for the remaining two data columns
$ws.Range("B1").EntireColumn.ColumnWidth = 23.8
$ws.Range("C1").EntireColumn.ColumnWidth = 25.8
